# Add a "Portugal" test-data sheet, cloned from the existing "Italy" sheet,
# and keep it as the active/selected tab (matching the target workbook diff).

$wb = $excel.ActiveWorkbook

# The Italy worksheet is the template for the new Portugal sheet.
$italy = $wb.Worksheets.Item("Italy")

# Duplicate Italy, placing the copy immediately after it; this becomes sheet #3
# ("Netherlands", "Italy", "Portugal") and is made the active sheet by Excel.
$italy.Copy($null, $italy)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Fill in the market name / user-story cells that differ from the Italy template.
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2551"

# Both the source sheet and the new sheet pick up an explicit portrait page setup.
$italy.PageSetup.Orientation = 1          # xlPortrait
$portugal.PageSetup.Orientation = 1       # xlPortrait

# Italy's old selection (A8:A15) is replaced by a "select all" of its grid, and
# it is no longer the selected tab now that Portugal has taken over - so select
# on Italy first, then hop back to Portugal so it stays the active sheet.
$italy.Activate()
$italy.Range("A1:XFD1048576").Select()

# Leave the cursor on B4 (the cell just edited) on the new, active Portugal sheet.
$portugal.Activate()
$portugal.Range("B4").Select()
